$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = 2
$ws.Range("C42").Value = 43
$ws.Range("D42").Value = 1
$ws.Range("E42").Value = "System"
$ws.Range("F42").Value = "2025-03-03 18:27:18"
$ws.Range("G42").Value = 0

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = 2
$ws.Range("C43").Value = 44
$ws.Range("D43").Value = 12
$ws.Range("E43").Value = "System"
$ws.Range("F43").Value = "2025-03-03 18:29:00"
$ws.Range("G43").Value = 0

$ws.Range("A44").Value = 43
$ws.Range("B44").Value = 2
$ws.Range("C44").Value = 45
$ws.Range("D44").Value = 23
$ws.Range("E44").Value = "System"
$ws.Range("F44").Value = "2025-03-03 18:31:04"
$ws.Range("G44").Value = 0
